$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Euler's method parameters:
#   dy/dt = 2*y - 2*t - 1 , h = 0.25, t0 = 0, y0 = 2
# Recompute rows 2..51 (t in column A, y in column B)

$h = 0.25
$t = 0.0
$y = 2.0

for ($i = 2; $i -le 51; $i++) {
    $ws.Cells.Item($i, 1).Value = $t
    $ws.Cells.Item($i, 2).Value = $y

    $f = 2 * $y - 2 * $t - 1
    $y = $y + $h * $f
    $t = $t + $h
}
